# Update parent column feature:
#  - Insert a new header sub-row ("Generation") under the existing
#    "Generations"/"Age"/"Birth" header block, rename the old
#    "Generations" header to "Life", widen the Birth/Category detail
#    column, and move the "Cartoon" category values into column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above the current row 3 ("Age"/"Birth" row); this
#    pushes the data rows down by one and Excel auto-extends the
#    existing A2:A3 / D2:D3 merges to A2:A4 / D2:D4.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).RowHeight = 30

# 2) Column widths: columns A, B, D stay 20.7109375; column C (now the
#    wide "Birth"/detail column) widens to ~40.71 (double the normal
#    column width).
$ws.Columns.Item(3).ColumnWidth = 39.8

# 3) Row 2 header cells.
#    B2: "Generations" -> "Life", and drop the bold/red style down to
#    the plain bold style used elsewhere in the header (style of A2).
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").Value = "Life"

#    C2 was a duplicate "Generations" label - clear its text, it stays
#    merged away behind B2:C2.
$ws.Range("C2").ClearContents()

# 4) New row 3 ("Generation" sub-header), styled like the rest of the
#    plain header row.
$ws.Range("A2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = "Generation"

$ws.Range("A2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# 5) Move the "Cartoon" category values from column C to column D on
#    the two data rows (now rows 5 and 6 after the insert).
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("C5").ClearContents()
$ws.Range("D6").Value = $ws.Range("C6").Value2
$ws.Range("C6").ClearContents()

# 6) New merged header cells.
$ws.Range("B2:C2").Merge()
$ws.Range("B3:C3").Merge()
